$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 13545
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 13545
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = $null
$ws.Range("M46").Value = 40635
$ws.Range("N46").Value = -40873

$ws.Range("H53").Value = 38747.117
$ws.Range("I53").Value = 77067.08
$ws.Range("J53").Value = 427.15384
$ws.Range("K53").Value = 77067.08
$ws.Range("L53").Value = 427.15384
$ws.Range("M53").Value = -76430.08
$ws.Range("N53").Value = -1701.15384

$ws.Range("H60").Value = 13545
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 13545
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = $null
$ws.Range("M60").Value = 40635
$ws.Range("N60").Value = -41603

$ws.Range("H76").Value = 42977.8
$ws.Range("I76").Value = 44560.207
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 44560.207
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -44245.207
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 42977.8
$ws.Range("I79").Value = 44560.207
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 44560.207
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -43468.207
$ws.Range("N79").Value = -7184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 14144.75
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 14144.75
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = $null
$ws.Range("M44").Value = 14144.75
$ws.Range("N44").Value = -15120.75

$ws.Range("H122").Value = 1569.4642
$ws.Range("I122").Value = 1517.2
$ws.Range("K122").Value = 4551.6
$ws.Range("M122").Value = -2101.6

$ws.Range("H132").Value = 7404.579
$ws.Range("I132").Value = 8263.357
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 24790.071
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -22260.071
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 250
$ws.Range("I8").Value = 250
$ws.Range("K8").Value = 750
$ws.Range("M8").Value = -611

$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = $null
$ws.Range("M15").Value = 900
$ws.Range("N15").Value = -1180

$ws.Range("H26").Value = 50000150
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 100000250
$ws.Range("K26").Value = 150
$ws.Range("L26").Value = 300000750
$ws.Range("M26").Value = 138
$ws.Range("N26").Value = -300001326

$ws.Range("H69").Value = 1414
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1414
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = $null
$ws.Range("M69").Value = 4242
$ws.Range("N69").Value = -5864

$ws.Range("H72").Value = 1414
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1414
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = $null
$ws.Range("M72").Value = 12726
$ws.Range("N72").Value = -20838

$ws.Range("H104").Value = 1200
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 1200
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = $null
$ws.Range("M104").Value = 3600
$ws.Range("N104").Value = -8842

$ws.Range("H113").Value = 512.129
$ws.Range("J113").Value = 495.81818
$ws.Range("L113").Value = 1487.45454
$ws.Range("N113").Value = -5827.45454

$ws.Range("H131").Value = 3835597.5
$ws.Range("J131").Value = 6536905
$ws.Range("L131").Value = 19610715
$ws.Range("N131").Value = -19620795

$ws.Range("H132").Value = 1625
$ws.Range("I132").Value = 490
$ws.Range("J132").Value = 2255.5557
$ws.Range("K132").Value = 4410
$ws.Range("L132").Value = 20300.0013
$ws.Range("M132").Value = -1880
$ws.Range("N132").Value = -25360.0013

$ws.Range("H137").Value = 173644180
$ws.Range("I137").Value = 111111740
$ws.Range("J137").Value = 211163630
$ws.Range("K137").Value = 333335220
$ws.Range("L137").Value = 633490890
$ws.Range("M137").Value = -333330120
$ws.Range("N137").Value = -633501090

$ws.Range("H141").Value = 2122.2068
$ws.Range("I141").Value = 1980.174
$ws.Range("K141").Value = 5940.522
$ws.Range("M141").Value = -760.5219999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7833.3335
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 9500
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 9500
$ws.Range("M43").Value = -4349
$ws.Range("N43").Value = -9802

$ws.Range("H46").Value = 19548.666
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 19548.666
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = $null
$ws.Range("M46").Value = 19548.666
$ws.Range("N46").Value = -19860.666

$ws.Range("H52").Value = 16500
$ws.Range("J52").Value = 16500
$ws.Range("L52").Value = 16500
$ws.Range("N52").Value = -17018

$ws.Range("H126").Value = 6118.6665
$ws.Range("I126").Value = 7102.4
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 21307.2
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -18837.2
$ws.Range("N126").Value = -8540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20835450
$ws.Range("I7").Value = 1779.6666
$ws.Range("J7").Value = 55558230
$ws.Range("K7").Value = 1779.6666
$ws.Range("L7").Value = 55558230
$ws.Range("M7").Value = -1667.6666
$ws.Range("N7").Value = -55558454

$ws.Range("H93").Value = 6758557
$ws.Range("I93").Value = 6758557
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 6758557
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = -6757309

$ws.Range("H126").Value = 20835450
$ws.Range("I126").Value = 1779.6666
$ws.Range("J126").Value = 55558230
$ws.Range("K126").Value = 5338.9998
$ws.Range("L126").Value = 166674690
$ws.Range("M126").Value = -2868.9998
$ws.Range("N126").Value = -166679630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1168.4615
$ws.Range("I122").Value = 1128.5714
$ws.Range("J122").Value = 1215
$ws.Range("K122").Value = 3385.7142
$ws.Range("L122").Value = 3645
$ws.Range("M122").Value = -935.7142000000003
$ws.Range("N122").Value = -8545

$ws.Range("H126").Value = 1038.2778
$ws.Range("I126").Value = 679.26666
$ws.Range("J126").Value = 2833.3333
$ws.Range("K126").Value = 2037.79998
$ws.Range("L126").Value = 8499.999899999999
$ws.Range("M126").Value = 432.20002
$ws.Range("N126").Value = -13439.9999

$ws.Range("H132").Value = 2006.0625
$ws.Range("I132").Value = 1424.5
$ws.Range("K132").Value = 4273.5
$ws.Range("M132").Value = -1743.5
